# Split the "KHĐC (Tăng 3-5%)" header column into two separate columns:
#   - "KHBĐ"  (kept in the original column, text shortened)
#   - "KHĐC"  (new column inserted immediately to its right)
#
# This mirrors a manual Excel edit: right-click column M header -> Insert,
# then retype the two header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank column at M; everything from the old M onward
# (M,N,O,P) shifts right by one (-> N,O,P,Q). K and L keep their positions.
$ws.Columns("M").Insert()

# The old combined header text lived in L2 ("KHĐC                       (Tăng 3-5%)").
# Shorten it in place to "KHBĐ" ...
$ws.Range("L2").Value = "KHBĐ"

# ...and give the newly inserted column M its own header, "KHĐC".
$ws.Range("M2").Value = "KHĐC                      "

# Leave the selection where the author left it before saving.
$null = $ws.Range("P4").Select()
